$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data columns stay plain text so values like "1.00", "9.00",
# "67.631.39" and "  -0.47%  " round-trip exactly instead of being coerced
# into numbers/dates by Excel's usual type inference.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 2).Value = "Bitcoin"
$ws.Cells.Item(2, 3).Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Cells.Item(2, 4).Value = "67.631.39"
$ws.Cells.Item(2, 5).Value = "  -0.47%  "
$ws.Cells.Item(3, 2).Value = "Ethereum"
$ws.Cells.Item(3, 3).Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Cells.Item(3, 4).Value = "3.723.41"
$ws.Cells.Item(3, 5).Value = "  -2.36%  "
$ws.Cells.Item(4, 2).Value = "TetherUSD"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.20%  "
$ws.Cells.Item(5, 2).Value = "BNB"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Cells.Item(5, 4).Value = "590.91"
$ws.Cells.Item(5, 5).Value = "  -1.46%  "
$ws.Cells.Item(6, 2).Value = "Solana"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(6, 4).Value = "165.25"
$ws.Cells.Item(6, 5).Value = "  -2.33%  "
$ws.Cells.Item(7, 2).Value = "LidoStakedEther"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(7, 4).Value = "3.723.45"
$ws.Cells.Item(7, 5).Value = "  -2.34%  "
$ws.Cells.Item(8, 2).Value = "USDC"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  +0.00%  "
$ws.Cells.Item(9, 2).Value = "XRP"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(9, 4).Value = "0.518"
$ws.Cells.Item(9, 5).Value = "  -2.33%  "
$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).Value = "0.159"
$ws.Cells.Item(10, 5).Value = "  -3.80%  "
$ws.Cells.Item(11, 2).Value = "Toncoin"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(11, 4).Value = "6.49"
$ws.Cells.Item(11, 5).Value = "  -0.26%  "
$ws.Cells.Item(12, 2).Value = "Cardano"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(12, 4).Value = "0.450"
$ws.Cells.Item(12, 5).Value = "  -2.85%  "
$ws.Cells.Item(13, 2).Value = "ShibaInu"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(13, 4).Value = "0.0000261"
$ws.Cells.Item(13, 5).Value = "  -5.45%  "
$ws.Cells.Item(14, 2).Value = "Avalanche"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(14, 4).Value = "36.13"
$ws.Cells.Item(14, 5).Value = "  -2.40%  "
$ws.Cells.Item(15, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(15, 4).Value = "4.357.03"
$ws.Cells.Item(15, 5).Value = "  -2.13%  "
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "3.728.28"
$ws.Cells.Item(16, 5).Value = "  -2.36%  "
$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "67.662.65"
$ws.Cells.Item(17, 5).Value = "  -0.57%  "
$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).Value = "18.25"
$ws.Cells.Item(18, 5).Value = "  -1.66%  "
$ws.Cells.Item(19, 2).Value = "Polkadot"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(19, 4).Value = "7.01"
$ws.Cells.Item(19, 5).Value = "  -5.99%  "
$ws.Cells.Item(20, 2).Value = "TRON"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(20, 4).Value = "0.111"
$ws.Cells.Item(20, 5).Value = "  -0.51%  "
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(21, 4).Value = "10.68"
$ws.Cells.Item(21, 5).Value = "  -1.96%  "
$ws.Cells.Item(22, 2).Value = "BitcoinCash"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(22, 4).Value = "466.91"
$ws.Cells.Item(22, 5).Value = "  -0.81%  "
$ws.Cells.Item(23, 2).Value = "Polygon"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(23, 4).Value = "0.698"
$ws.Cells.Item(23, 5).Value = "  -5.69%  "
$ws.Cells.Item(24, 2).Value = "Litecoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(24, 4).Value = "82.69"
$ws.Cells.Item(24, 5).Value = "  -0.96%  "
$ws.Cells.Item(25, 2).Value = "Fetch.AI"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(25, 4).Value = "2.18"
$ws.Cells.Item(25, 5).Value = "  -5.45%  "
$ws.Cells.Item(26, 2).Value = "PEPE"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(26, 4).Value = "0.0000133"
$ws.Cells.Item(26, 5).Value = "  -12.59%  "
$ws.Cells.Item(27, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 4).Value = "12.01"
$ws.Cells.Item(27, 5).Value = "  -1.56%  "
$ws.Cells.Item(28, 2).Value = "RenderToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(28, 4).Value = "10.13"
$ws.Cells.Item(28, 5).Value = "  -1.33%  "
$ws.Cells.Item(29, 2).Value = "Dai"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 5).Value = "  +0.00%  "
$ws.Cells.Item(30, 2).Value = "WrappedeETH"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(30, 4).Value = "3.876.01"
$ws.Cells.Item(30, 5).Value = "  -2.08%  "
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).Value = "2.76"
$ws.Cells.Item(31, 5).Value = "  -5.60%  "
$ws.Cells.Item(32, 2).Value = "NEARProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(32, 4).Value = "7.30"
$ws.Cells.Item(32, 5).Value = "  -5.59%  "
$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(33, 4).Value = "2.21"
$ws.Cells.Item(33, 5).Value = "  -3.38%  "
$ws.Cells.Item(34, 2).Value = "EthereumClassic"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(34, 4).Value = "29.48"
$ws.Cells.Item(34, 5).Value = "  -4.26%  "
$ws.Cells.Item(35, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(35, 4).Value = "0.997"
$ws.Cells.Item(35, 5).Value = "  --%  "
$ws.Cells.Item(36, 2).Value = "Aptos"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(36, 4).Value = "9.00"
$ws.Cells.Item(36, 5).Value = "  -3.91%  "
$ws.Cells.Item(37, 2).Value = "RenzoRestakedETH"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Cells.Item(37, 4).Value = "3.682.69"
$ws.Cells.Item(37, 5).Value = "  -2.54%  "
$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).Value = "0.101"
$ws.Cells.Item(38, 5).Value = "  -5.79%  "
$ws.Cells.Item(39, 2).Value = "dogwifhat"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(39, 4).Value = "3.39"
$ws.Cells.Item(39, 5).Value = "  -11.21%  "
$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).Value = "0.138"
$ws.Cells.Item(40, 5).Value = "  -1.48%  "
$ws.Cells.Item(41, 2).Value = "Mantle"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(41, 4).Value = "0.987"
$ws.Cells.Item(41, 5).Value = "  -2.33%  "
$ws.Cells.Item(42, 2).Value = "Filecoin"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(42, 4).Value = "5.73"
$ws.Cells.Item(42, 5).Value = "  -4.60%  "
$ws.Cells.Item(43, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(43, 4).Value = "0.999"
$ws.Cells.Item(43, 5).Value = "  +0.16%  "
$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(44, 4).Value = "1.00"
$ws.Cells.Item(44, 5).Value = "  -0.01%  "
$ws.Cells.Item(45, 2).Value = "TheGraph"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(45, 4).Value = "0.304"
$ws.Cells.Item(45, 5).Value = "  -4.70%  "
$ws.Cells.Item(46, 2).Value = "Cosmos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(46, 4).Value = "8.51"
$ws.Cells.Item(46, 5).Value = "  -3.47%  "
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).Value = "1.91"
$ws.Cells.Item(47, 5).Value = "  -3.43%  "
$ws.Cells.Item(48, 2).Value = "OKB"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(48, 4).Value = "45.25"
$ws.Cells.Item(48, 5).Value = "  -2.57%  "
$ws.Cells.Item(49, 2).Value = "Bittensor"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(49, 4).Value = "386.58"
$ws.Cells.Item(49, 5).Value = "  -6.01%  "
$ws.Cells.Item(50, 2).Value = "Monero"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(50, 4).Value = "143.00"
$ws.Cells.Item(50, 5).Value = "  +0.74%  "
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "25.24"
$ws.Cells.Item(51, 5).Value = "  -0.23%  "
